$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O (15) to make room for "Birth Year"
$ws.Columns.Item(15).Insert()

# Set the header for the new column
$ws.Cells.Item(1, 15).Value = "Birth Year"

# New column inherits the same width as "Age Units" (column N) so the two
# columns render the same as in the published template
$ws.Columns.Item(15).ColumnWidth = $ws.Columns.Item(14).ColumnWidth

# Scroll / selection adjustments to match the saved view state
$ws.Range("M2").Select()
$excel.ActiveWindow.ScrollColumn = 8
